$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 1
$ws.Range("F24").Value = 0
$ws.Range("H24").Value = 1
$ws.Range("F46").Value = 0
$ws.Range("H46").Value = 1
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("E65").Value = 2
$ws.Range("H65").Value = 4
$ws.Range("D66").Value = 2
$ws.Range("H66").Value = 4
$ws.Range("F67").Value = 2
$ws.Range("H67").Value = 3
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 3
$ws.Range("F70").Value = 0
$ws.Range("H70").Value = 1
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("E86").Value = 2
$ws.Range("H86").Value = 3
$ws.Range("F87").Value = 1
$ws.Range("H87").Value = 1
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 1
$ws.Range("D97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 0
$ws.Range("D101").Value = 1
$ws.Range("H101").Value = 1
$ws.Range("E102").Value = 1
$ws.Range("H102").Value = 1
$ws.Range("F103").Value = 1
$ws.Range("H103").Value = 1
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 1
$ws.Range("D141").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("F142").Value = 0
$ws.Range("H142").Value = 1
$ws.Range("E143").Value = 0
$ws.Range("H143").Value = 1
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0
$ws.Range("E149").Value = 1
$ws.Range("H149").Value = 3
$ws.Range("F159").Value = 1
$ws.Range("H159").Value = 2
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 0
$ws.Range("E162").Value = 1
$ws.Range("H162").Value = 2
$ws.Range("F163").Value = 1
$ws.Range("H163").Value = 2
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 0
$ws.Range("F175").Value = 0
$ws.Range("H175").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0
